# tasks.xlsx - "Auto update md files" edit
#
# Semantic changes applied (see xml diff):
#   1. F2 (remarks for "完成下周工作计划"): rich-text note "需要在周四下午4.前完成"
#      is replaced by the short plain-text note "急".
#   2. E6 (status for "超市采购生活用品"): dropdown value changes from "已完成"
#      to "未完成" (its G6 system-judged status cell recalculates accordingly).
#   3. B7 (due date for "整理书房电子发票"): due date moves from 46052 (2026-01-30)
#      to 46049 (2026-01-27); G7 recalculates accordingly.
#   4. Column F is widened (remarks column grew to fit the longer due-date note).
#   5. The saved cursor/selection position moves to F17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shrink the "备注" (remarks) note on row 2 down to a short plain-text "急" (urgent).
$ws.Range("F2").Value = "急"

# 2. Flip the "当前状态" (current status) of the grocery-shopping task to "未完成".
$ws.Range("E6").Value = "未完成"

# 3. Push back the due date of the "整理书房电子发票" task.
$ws.Range("B7").Value = 46049

# 4. Widen the remarks column (F) to match the new content width.
$ws.Columns.Item(6).ColumnWidth = 45.5714285714286

# 5. Leave the cursor where the author last left it.
$ws.Range("F17").Select()
